$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Structural changes: make room for the new "status_id"/"series_id" columns
# on the Product table, and for two new small lookup tables (Status,
# Category) plus the new Series table appended below OrderItem.
# ---------------------------------------------------------------------------

# Insert two new columns (K:L) into the Product header row for status_id /
# series_id; this pushes the Color/Storage/OrderItem side-tables from
# columns L:N out to N:P.
$ws.Columns("K:L").Insert()

# Insert a row to push the Order/Storage tables down from row 7/8 to row 8/9
# (inserted a couple of rows below the Product header so Excel doesn't carry
# the header formatting down into the new row).
$ws.Rows("6:6").Insert()

# Insert two rows to push the User/OrderItem tables down from row 11/12 to
# row 14/15, leaving room for the new Category table in between.
$ws.Rows("12:13").Insert()

# ---------------------------------------------------------------------------
# Values, written in the same order the table was authored so newly
# interned shared strings line up with the source workbook.
# ---------------------------------------------------------------------------

# Product table: two new header columns
$ws.Range("K3").Value = "status_id"

# New Status table (title + header)
$ws.Range("O5").Value = "Status"

# Color table gains a "value" column
$ws.Range("P3").Value = "value"

# New Category table (title + header)
$ws.Range("O11").Value = "Category"
$ws.Range("N12").Value = "category_id"
$ws.Range("O12").Value = "category"

# New Series table (title + header), appended below OrderItem
$ws.Range("O18").Value = "Series"
$ws.Range("N19").Value = "series_id"
$ws.Range("P19").Value = "series"

# Remaining cells reuse strings already interned above
$ws.Range("L3").Value = "series_id"
$ws.Range("N6").Value = "status_id"
$ws.Range("O6").Value = "status"
$ws.Range("P6").Value = "value"
$ws.Range("P9").Value = "value"
$ws.Range("O19").Value = "category_id"

# ---------------------------------------------------------------------------
# Formatting: copy the look of sibling title/header cells onto the new ones.
# ---------------------------------------------------------------------------
$ws.Range("N3").Copy()
$ws.Range("P3").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("O2").Copy()
$ws.Range("O5").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("O11").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("O18").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("N3").Copy()
$ws.Range("N6:P6").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("N12:O12").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("N19:P19").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("N9").Copy()
$ws.Range("P9").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("G4").Select()
